# Auto-generated Excel COM-interop script
# Applies numeric corrections to the Leve-profit tracking sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled price-refresh run.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 3457072.8
$ws.Range("I15").Value = 3457072.8
$ws.Range("K15").Value = 10371218.4
$ws.Range("M15").Value = -10371049.4
# Row 20
$ws.Range("H20").Value = 10010.5
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
# Row 21
$ws.Range("H21").Value = 10008.5
$ws.Range("J21").Value = 15000
$ws.Range("L21").Value = 15000
$ws.Range("N21").Value = -15936
# Row 23
$ws.Range("H23").Value = 10008.5
$ws.Range("J23").Value = 15000
$ws.Range("L23").Value = 15000
$ws.Range("N23").Value = -15468
# Row 34
$ws.Range("H34").Value = 4732.7144
$ws.Range("I34").Value = 3855
$ws.Range("K34").Value = 3855
$ws.Range("M34").Value = -3652
# Row 35
$ws.Range("H35").Value = 10010.5
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
# Row 36
$ws.Range("H36").Value = 4732.7144
$ws.Range("I36").Value = 3855
$ws.Range("K36").Value = 3855
$ws.Range("M36").Value = -3140
# Row 39
$ws.Range("H39").Value = 438.3125
$ws.Range("I39").Value = 193.3077
$ws.Range("K39").Value = 579.9231
$ws.Range("M39").Value = -283.9231
# Row 40
$ws.Range("H40").Value = 3699
$ws.Range("I40").Value = 4875.5
$ws.Range("J40").Value = 3271.182
$ws.Range("K40").Value = 4875.5
$ws.Range("L40").Value = 3271.182
$ws.Range("M40").Value = -4700.5
$ws.Range("N40").Value = -3621.182
# Row 42
$ws.Range("H42").Value = 601.2632
$ws.Range("I42").Value = 636.4286
$ws.Range("J42").Value = 502.8
$ws.Range("K42").Value = 1909.2858
$ws.Range("L42").Value = 1508.4
$ws.Range("M42").Value = -1679.2858
$ws.Range("N42").Value = -1968.4
# Row 43
$ws.Range("H43").Value = 13706.571
$ws.Range("J43").Value = 7059.2
$ws.Range("L43").Value = 7059.2
$ws.Range("N43").Value = -7197.2
# Row 111
$ws.Range("H111").Value = 1731.4667
$ws.Range("J111").Value = 1802
$ws.Range("L111").Value = 5406
$ws.Range("N111").Value = -11540
# Row 116
$ws.Range("H116").Value = 4469
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 4469
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 4469
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -11353
# Row 132
$ws.Range("H132").Value = 1227.1666
$ws.Range("I132").Value = 1029.4
$ws.Range("J132").Value = 2216
$ws.Range("K132").Value = 3088.2
$ws.Range("L132").Value = 6648
$ws.Range("M132").Value = -558.2000000000003
$ws.Range("N132").Value = -11708
# Row 134
$ws.Range("H134").Value = 120000
$ws.Range("J134").Value = 120000
$ws.Range("L134").Value = 120000
$ws.Range("N134").Value = -130140
# Row 137
$ws.Range("H137").Value = 5111655
$ws.Range("I137").Value = 10003695
$ws.Range("K137").Value = 30011085
$ws.Range("M137").Value = -30008535
# Row 138
$ws.Range("H138").Value = 3302.878
$ws.Range("J138").Value = 3621.6
$ws.Range("L138").Value = 10864.8
$ws.Range("N138").Value = -21144.8

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1062.902
$ws.Range("I32").Value = 1062.902
$ws.Range("K32").Value = 1062.902
$ws.Range("M32").Value = -775.902
# Row 61
$ws.Range("H61").Value = 4612
$ws.Range("I61").Value = 4121.16
$ws.Range("K61").Value = 4121.16
$ws.Range("M61").Value = -3909.16
# Row 122
$ws.Range("H122").Value = 2628.5
$ws.Range("I122").Value = 2628.5
$ws.Range("K122").Value = 7885.5
$ws.Range("M122").Value = -5435.5
# Row 136
$ws.Range("H136").Value = 4612
$ws.Range("I136").Value = 4121.16
$ws.Range("K136").Value = 12363.48
$ws.Range("M136").Value = -9813.48

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 1506.619
$ws.Range("I99").Value = 1282
$ws.Range("J99").Value = 5999
$ws.Range("K99").Value = 1282
$ws.Range("L99").Value = 5999
$ws.Range("M99").Value = 216
$ws.Range("N99").Value = -8995
# Row 105
$ws.Range("H105").Value = 55572284
$ws.Range("I105").Value = 66685924
$ws.Range("J105").Value = 4066.6667
$ws.Range("K105").Value = 66685924
$ws.Range("L105").Value = 4066.6667
$ws.Range("M105").Value = -66684177
$ws.Range("N105").Value = -7560.6667

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 20835606
$ws.Range("I31").Value = 30304034
$ws.Range("K31").Value = 30304034
$ws.Range("M31").Value = -30303739
# Row 34
$ws.Range("H34").Value = 20835606
$ws.Range("I34").Value = 30304034
$ws.Range("K34").Value = 30304034
$ws.Range("M34").Value = -30303832
# Row 44
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 8551495
$ws.Range("I131").Value = 19608574
$ws.Range("J131").Value = 7389.227
$ws.Range("K131").Value = 58825722
$ws.Range("L131").Value = 22167.681
$ws.Range("M131").Value = -58820682
$ws.Range("N131").Value = -32247.681

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 48
$ws.Range("H48").Value = 42481.332
$ws.Range("J48").Value = 42481.332
$ws.Range("L48").Value = 42481.332
$ws.Range("N48").Value = -43451.332
# Row 53
$ws.Range("H53").Value = 59650.25
$ws.Range("J53").Value = 59770.332
$ws.Range("L53").Value = 59770.332
$ws.Range("N53").Value = -61032.332
# Row 136
$ws.Range("H136").Value = 54978.8
$ws.Range("J136").Value = 54978.8
$ws.Range("L136").Value = 164936.4
$ws.Range("N136").Value = -170036.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 54154.156
$ws.Range("I7").Value = 60319.41
$ws.Range("K7").Value = 60319.41
$ws.Range("M7").Value = -60207.41
# Row 22
$ws.Range("H22").Value = 3096.158
$ws.Range("I22").Value = 1777.4
$ws.Range("J22").Value = 4561.4443
$ws.Range("K22").Value = 1777.4
$ws.Range("L22").Value = 4561.4443
$ws.Range("M22").Value = -1482.4
$ws.Range("N22").Value = -5151.4443
# Row 27
$ws.Range("H27").Value = 3096.158
$ws.Range("I27").Value = 1777.4
$ws.Range("J27").Value = 4561.4443
$ws.Range("K27").Value = 1777.4
$ws.Range("L27").Value = 4561.4443
$ws.Range("M27").Value = -1670.4
$ws.Range("N27").Value = -4775.4443
# Row 55
$ws.Range("H55").Value = 624
$ws.Range("I55").Value = 825.25
$ws.Range("J55").Value = 355.66666
$ws.Range("K55").Value = 825.25
$ws.Range("L55").Value = 355.66666
$ws.Range("M55").Value = -652.25
$ws.Range("N55").Value = -701.66666
# Row 126
$ws.Range("H126").Value = 54154.156
$ws.Range("I126").Value = 60319.41
$ws.Range("K126").Value = 180958.23
$ws.Range("M126").Value = -178488.23

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 1781.9412
$ws.Range("I107").Value = 1349
$ws.Range("K107").Value = 4047
$ws.Range("M107").Value = -2127
# Row 132
$ws.Range("H132").Value = 5995.2
$ws.Range("I132").Value = 3651.4285
$ws.Range("K132").Value = 10954.2855
$ws.Range("M132").Value = -8424.2855

Write-Host "Applied scheduled profit-sheet corrections."
